$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("策略更新")

# Row 3: change exchange to 大连, day/night to 日盘, update date, content to 锁仓测试, status stays passed
$ws.Range("B3").Value = "大连"
$ws.Range("C3").Value = "日盘"
$ws.Range("D3").Value = Get-Date -Year 2017 -Month 2 -Day 21
$ws.Range("E3").Value = "锁仓测试"
$ws.Range("F3").Value = "passed"

# Row 4: clear all content (B4:F4), keep A4 as-is (already empty)
$ws.Range("B4:F4").ClearContents()

$wb.Save()
